$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Start the new "TLDR" corpus row in the columns that sit to the left
#    of the upcoming column insertion (A/B are untouched by inserting at
#    C), so the "TLDR" shared string is created first - same order the
#    original author's shared-string table shows.
$ws.Range("A7").Value = "TLDR"
$ws.Range("B7").Value = "EN"

# 2) Insert a new column before the old "Referenz" column (C) to hold
#    "Branche". This shifts the old C/D/E (Referenz/Größe/Anzahl) into
#    D/E/F, matching the diff.
$ws.Columns("C").Insert()

# 3) Fill the new "Branche" column top-to-bottom, including the header.
$ws.Range("C1").Value = "Branche"
$ws.Range("C2").Value = "Wiki"
$ws.Range("C3").Value = "Recht"
$ws.Range("C4").Value = "Wiki"
$ws.Range("C5").Value = "News"
$ws.Range("C6").Value = "News"
$ws.Range("C7").Value = "Social"

# 4) Finish the new row's remaining (post-insert) columns. The "Größe"
#    column (E) is right-aligned, matching the rest of that column.
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 11.8
$ws.Range("E7").HorizontalAlignment = -4152
$ws.Range("F7").Value = 4000000
